$wb = $excel.ActiveWorkbook

# --- Sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename the second sheet (first sheet name stays the same: "_set_FLOWS")
$ws2.Name = "_set_FLOWS_AGGREGATED"

# --- Sheet 1 ("_set_FLOWS") header rename ---
$ws1.Range("A1").Value = "flows_Name"
$ws1.Range("B1").Value = "flows_dispatch"
$ws1.Range("C1").Value = "flows_Aggregation"

# Data rows 2-6 are unchanged in content.

# --- Sheet 2 ("_set_FLOWS_AGGREGATED") header rename + new column ---
$ws2.Range("A1").Value = "flows_aggregated_Name"
$ws2.Range("B1").Value = "flows_aggregated_dispatch"
$ws2.Range("C1").Value = "flows_aggregated_Aggregation"
# Match the header formatting (bold/border/center) used by A1:B1
$ws2.Range("A1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-4 are unchanged in content (columns A and B only).

# --- Selections / active sheet state ---
# Sheet2 becomes the non-active sheet with selection A2:B4
$ws2.Activate()
$ws2.Range("A2:B4").Select()

# Sheet1 becomes the active/selected sheet with selection A2:C6
$ws1.Activate()
$ws1.Range("A2:C6").Select()
